$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("boulders")

# Insert 3 blank rows after row 14 ("above the ranch house") for new boulders,
# then 1 more blank row after what becomes row 18 ("sign of the cross").
$ws.Rows("15:17").Insert() | Out-Null
$ws.Rows("19:19").Insert() | Out-Null

# Reference formatting for the checkbox column (C) used by existing boulder rows.
$checkColor = $ws.Range("C2").Interior.Color

# Row 14 - above the ranch house: add new boulder "The Kitchen"
$ws.Range("B14").Value = "The Kitchen"
$ws.Range("C14").Interior.Color = $checkColor

# Row 15 (new) - "Power of Silence"
$ws.Range("B15").Value = "Power of Silence"
$ws.Range("C15").Interior.Color = $checkColor

# Row 16 (new) - "Gums Boulder"
$ws.Range("B16").Value = "Gums Boulder"
$ws.Range("C16").Interior.Color = $checkColor

# Row 17 (new) - "Look Sharp Rock"
$ws.Range("B17").Value = "Look Sharp Rock"
$ws.Range("C17").Interior.Color = $checkColor

# Row 18 - sign of the cross: add new boulder "Sign of the cross wall"
$ws.Range("B18").Value = "Sign of the cross wall"
$ws.Range("C18").Interior.Color = $checkColor

# Row 19 (new) - "Term Boulder"
$ws.Range("B19").Value = "Term Boulder"
$ws.Range("C19").Interior.Color = $checkColor

# Update selection to reflect the author's last edit position.
$ws.Range("C18:C19").Select() | Out-Null
